$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.89490282535553
$ws.Range("B1").Value = 2.08281421661377
$ws.Range("C1").Value = 2.458517074584961
$ws.Range("D1").Value = 3.072900295257568
$ws.Range("E1").Value = 2.386813879013062
